# Appends the new bullet-list items documenting the "Natours part 2" section
# (Tours / Stories / Booking / Navigation sections and the CSS-only popup)
# right after the existing "Feature section: ..." bullet, before the final
# sectPr. Each new paragraph reuses the ListParagraph style / numId=12 bullet
# numbering and Roboto 12pt (sz 24 half-points) run formatting already used by
# the rest of that bulleted list, because InsertParagraphAfter() on the last
# paragraph of the document copies that paragraph's pPr/rPr.

$d = $word.ActiveDocument

function Add-BulletParagraph($RunTexts) {
    $anchor = $d.Paragraphs.Last.Range
    $anchor.Collapse(0)
    $anchor.InsertParagraphAfter()

    $newRange = $d.Paragraphs.Last.Range
    $newRange.Collapse(0)
    foreach ($t in $RunTexts) {
        $newRange.InsertAfter($t)
        $newRange.Collapse(0)
    }
}

# Tours section
Add-BulletParagraph @('Tours section: How to build an amazing, rotating card; -How to use perspective in CSS; 0How to use the backface-visibility property; -Using background blend modes; -How and when to use box-decoration-break.')

# the box-decoration-break explanation paragraph
Add-BulletParagraph @(
    '//',
    'if the text is split on 2 lines we can interpret each of them as two separate elements',
    ' and we can apply the decorations on each of them',
    ' ->',
    ' box-decoration-break: clone;'
)

# Stories section
Add-BulletParagraph @(
    'Stories section: How to make text flow around shapes with shape-outside and float; -how to apply a filter to images; -hot to create a background video covering an entire section; -how to use the <video> HTML element; -how and when to use the object-fit property',
    '(this is to cover the parent with the html element – in my case with a background video)'
)

# Booking section
Add-BulletParagraph @('Booking section: -How to implement “solid-color gradients”; -How the generatlr and adjacent sibling celectors work and why we need them; -How to use the::input-placeholder pseudo-element; -how and when to use the :focus, :invalid, placeholder-shown and :checked pseudo-classes; -Techniques to build custom radio buttons;')

# Navigation part
Add-BulletParagraph @('Navigation part: What the “checkbox hack” is and how it works; -How to create custom animation timing functions using cubic Bezier curves; -How to animate “solid-color gradients”; -How and why to use transform-origin; -In generatl: create an amazingly creative effect;')

# Building a popup with only CSS
Add-BulletParagraph @(
    'Building a popup with only CSS: -How to build a nice popup with only CSS ; -How to use the :target pseudo-class; -how to create boxes with equal height using ',
    'display: table-cell; -how to create CSS text columns; -How to automatically hyphenate words using hypens'
)
